$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("J3").Value = 1.14
$ws.Range("K3").Formula = "=0.48/2"

# Row 5
$ws.Range("J5").Value = 1.23

# Row 8
$ws.Range("J8").Value = 1.04
$ws.Range("K8").Formula = "=0.61/2"

# Row 9
$ws.Range("J9").Value = 0.85
$ws.Range("K9").Formula = "=0.49/2"

# Row 15
$ws.Range("J15").Value = 1.09
$ws.Range("K15").Formula = "=0.93/2"

# Update selected cell to K16
$ws.Range("K16").Select()
